$d = $word.ActiveDocument

$replacements = @(
    @("31×88=2728", "92×80=7360"),
    @("24×75=1800", "19×99=1881"),
    @("91×71=6461", "31×77=2387"),
    @("15×87=1305", "45×24=1080"),
    @("96×19=1824", "78×31=2418"),
    @("64×58=3712", "13×50=650"),
    @("53×91=4823", "83×97=8051"),
    @("58×31=1798", "86×79=6794"),
    @("81×96=7776", "88×62=5456"),
    @("54×53=2862", "64×55=3520"),
    @("65×95=6175", "79×24=1896"),
    @("51×61=3111", "24×94=2256"),
    @("55×37=2035", "89×89=7921"),
    @("96×90=8640", "20×47=940"),
    @("85×51=4335", "18×71=1278"),
    @("43×45=1935", "13×78=1014"),
    @("38×57=2166", "98×14=1372"),
    @("50×85=4250", "64×37=2368"),
    @("92×54=4968", "36×59=2124"),
    @("54×75=4050", "34×81=2754"),
    @("41×46=1886", "67×43=2881"),
    @("51×66=3366", "45×70=3150"),
    @("34×85=2890", "91×37=3367"),
    @("16×43=688", "46×88=4048"),
    @("40×12=480", "31×11=341")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
